# Fruta / hortaliza, semanal
# The weekly data refresh reshuffles each price-record (columns D, M-T)
# among the existing rows 2-26. Columns A,B,C,E-L (market/product metadata)
# are identical on every row, so only D,M,N,O,P,Q,R,S,T need to move.
#
# Mapping: new row -> source row that its record comes from.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 6
    3  = 5
    4  = 11
    5  = 21
    6  = 18
    7  = 4
    8  = 7
    9  = 16
    10 = 10
    11 = 17
    12 = 22
    13 = 12
    14 = 8
    15 = 3
    16 = 14
    17 = 26
    18 = 19
    19 = 25
    20 = 13
    21 = 23
    22 = 2
    23 = 24
    24 = 9
    25 = 20
    26 = 15
}

# Columns that move as a record (1-based column indices):
# D=4, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20
$cols = @(4, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot current values for every row/column before writing anything,
# since this is a permutation and rows read-from and written-to overlap.
# NOTE: use Value2 (not Value) when reading — Value round-trips cleanly here.
$snapshot = @{}
for ($r = 2; $r -le 26; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the permuted values to their new homes.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $srcVals = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value2 = $srcVals[$c]
    }
}
